$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.317.53"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").Value = "3.488.41"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.88"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +12.47%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "3.489.77"
$ws.Range("E9").Value = "  +2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("E12").Value = "  +3.83%  "
$ws.Range("D13").Value = "4.090.51"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000194"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.78"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.76%  "
$ws.Range("D17").Value = "65.299.79"
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").Value = "3.485.28"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.47"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "384.31"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.24"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.555"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.77%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.73"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000121"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.27%  "
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("E29").Value = "  +12.90%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.76"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.22"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.05%  "
$ws.Range("E35").Value = "  +13.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.19"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.94"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +6.73%  "
$ws.Range("D38").Value = "3.015.44"
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0783"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.54%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.98"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.81"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.29%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.59"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.98%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0324"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.97"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.784"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.09"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +11.90%  "
$ws.Range("E47").Value = "  +4.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "320.44"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +9.40%  "
$ws.Range("E49").Value = "  +6.88%  "
$ws.Range("E50").Value = "  +5.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.880"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.28%  "
